# Auto-generated edit script applying updated profit-calculation values
# across multiple sheets/rows in the workbook (columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 98259.664
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 98259.664
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 98259.664
$ws.Range("N123").Value = -108059.664

$ws.Range("H137").Value = 1603.6522
$ws.Range("I137").Value = 1603.6522
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4810.9566
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2260.9566
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 1772.98
$ws.Range("I138").Value = 877.3684
$ws.Range("J138").Value = 1983.0618
$ws.Range("K138").Value = 2632.1052
$ws.Range("L138").Value = 5949.1854
$ws.Range("M138").Value = 2507.8948
$ws.Range("N138").Value = -16229.1854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2821.55
$ws.Range("I61").Value = 2033.3334
$ws.Range("J61").Value = 4003.875
$ws.Range("K61").Value = 2033.3334
$ws.Range("L61").Value = 4003.875
$ws.Range("M61").Value = -1821.3334
$ws.Range("N61").Value = -4427.875

$ws.Range("H63").Value = 9789.5
$ws.Range("I63").Value = 10141.471
$ws.Range("J63").Value = 3806
$ws.Range("K63").Value = 10141.471
$ws.Range("L63").Value = 3806
$ws.Range("M63").Value = -9455.471
$ws.Range("N63").Value = -5178

$ws.Range("H66").Value = 9789.5
$ws.Range("I66").Value = 10141.471
$ws.Range("J66").Value = 3806
$ws.Range("K66").Value = 50707.355
$ws.Range("L66").Value = 19030
$ws.Range("M66").Value = -47275.355
$ws.Range("N66").Value = -25894

$ws.Range("H74").Value = 6677.4546
$ws.Range("I74").Value = 1304.7368
$ws.Range("J74").Value = 40704.668
$ws.Range("K74").Value = 1304.7368
$ws.Range("L74").Value = 40704.668
$ws.Range("M74").Value = -430.7367999999999
$ws.Range("N74").Value = -42452.668

$ws.Range("H77").Value = 6677.4546
$ws.Range("I77").Value = 1304.7368
$ws.Range("J77").Value = 40704.668
$ws.Range("K77").Value = 6523.683999999999
$ws.Range("L77").Value = 203523.34
$ws.Range("M77").Value = -2155.683999999999
$ws.Range("N77").Value = -212259.34

$ws.Range("H122").Value = 7945.4
$ws.Range("I122").Value = 10029
$ws.Range("J122").Value = 2215.5
$ws.Range("K122").Value = 30087
$ws.Range("L122").Value = 6646.5
$ws.Range("M122").Value = -27637
$ws.Range("N122").Value = -11546.5

$ws.Range("H123").Value = 32523.834
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 32523.834
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 32523.834
$ws.Range("N123").Value = -42323.834

$ws.Range("H136").Value = 2821.55
$ws.Range("I136").Value = 2033.3334
$ws.Range("J136").Value = 4003.875
$ws.Range("K136").Value = 6100.0002
$ws.Range("L136").Value = 12011.625
$ws.Range("M136").Value = -3550.0002
$ws.Range("N136").Value = -17111.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 35000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 35000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040

$ws.Range("H134").Value = 4191.8887
$ws.Range("I134").Value = 3311.2
$ws.Range("J134").Value = 5292.75
$ws.Range("K134").Value = 9933.599999999999
$ws.Range("L134").Value = 15878.25
$ws.Range("M134").Value = -7398.599999999999
$ws.Range("N134").Value = -20948.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5093.4644
$ws.Range("I31").Value = 1998.4286
$ws.Range("J31").Value = 6125.143
$ws.Range("K31").Value = 1998.4286
$ws.Range("L31").Value = 6125.143
$ws.Range("M31").Value = -1703.4286
$ws.Range("N31").Value = -6715.143

$ws.Range("H34").Value = 5093.4644
$ws.Range("I34").Value = 1998.4286
$ws.Range("J34").Value = 6125.143
$ws.Range("K34").Value = 1998.4286
$ws.Range("L34").Value = 6125.143
$ws.Range("M34").Value = -1796.4286
$ws.Range("N34").Value = -6529.143

$ws.Range("H58").Value = 2396.75
$ws.Range("I58").Value = 1217.3636
$ws.Range("J58").Value = 3838.2222
$ws.Range("K58").Value = 1217.3636
$ws.Range("L58").Value = 3838.2222
$ws.Range("M58").Value = -1014.3636
$ws.Range("N58").Value = -4244.2222

$ws.Range("H134").Value = 4144.0454
$ws.Range("I134").Value = 2598.6428
$ws.Range("J134").Value = 6848.5
$ws.Range("K134").Value = 7795.928400000001
$ws.Range("L134").Value = 20545.5
$ws.Range("M134").Value = -5260.928400000001
$ws.Range("N134").Value = -25615.5

$ws.Range("H136").Value = 2396.75
$ws.Range("I136").Value = 1217.3636
$ws.Range("J136").Value = 3838.2222
$ws.Range("K136").Value = 3652.0908
$ws.Range("L136").Value = 11514.6666
$ws.Range("M136").Value = -1102.0908
$ws.Range("N136").Value = -16614.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2247.25
$ws.Range("I5").Value = 1972.6666
$ws.Range("J5").Value = 2600.2856
$ws.Range("K5").Value = 5917.9998
$ws.Range("L5").Value = 7800.8568
$ws.Range("M5").Value = -5805.9998
$ws.Range("N5").Value = -8024.8568

$ws.Range("H135").Value = 2247.25
$ws.Range("I135").Value = 1972.6666
$ws.Range("J135").Value = 2600.2856
$ws.Range("K135").Value = 17753.9994
$ws.Range("L135").Value = 23402.5704
$ws.Range("M135").Value = -15218.9994
$ws.Range("N135").Value = -28472.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6895.4
$ws.Range("I70").Value = 8262.154
$ws.Range("J70").Value = 4357.143
$ws.Range("K70").Value = 8262.154
$ws.Range("L70").Value = 4357.143
$ws.Range("M70").Value = -7992.154
$ws.Range("N70").Value = -4897.143

$ws.Range("H73").Value = 6895.4
$ws.Range("I73").Value = 8262.154
$ws.Range("J73").Value = 4357.143
$ws.Range("K73").Value = 8262.154
$ws.Range("L73").Value = 4357.143
$ws.Range("M73").Value = -7326.154
$ws.Range("N73").Value = -6229.143

$ws.Range("H102").Value = 2225.5
$ws.Range("I102").Value = 1496
$ws.Range("J102").Value = 2468.6667
$ws.Range("K102").Value = 1496
$ws.Range("L102").Value = 2468.6667
$ws.Range("M102").Value = 126
$ws.Range("N102").Value = -5712.6667

$ws.Range("H107").Value = 845.6
$ws.Range("I107").Value = 1123.7778
$ws.Range("J107").Value = 428.33334
$ws.Range("K107").Value = 1123.7778
$ws.Range("L107").Value = 428.33334
$ws.Range("M107").Value = 796.2221999999999
$ws.Range("N107").Value = -4268.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5032
$ws.Range("I61").Value = 6450.75
$ws.Range("J61").Value = 2762
$ws.Range("K61").Value = 6450.75
$ws.Range("L61").Value = 2762
$ws.Range("M61").Value = -6248.75
$ws.Range("N61").Value = -3166

$ws.Range("H94").Value = 24999.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 24999.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 24999.5
$ws.Range("N94").Value = -26351.5

$ws.Range("H113").Value = 5032
$ws.Range("I113").Value = 6450.75
$ws.Range("J113").Value = 2762
$ws.Range("K113").Value = 6450.75
$ws.Range("L113").Value = 2762
$ws.Range("M113").Value = -4280.75
$ws.Range("N113").Value = -7102

$ws.Range("H132").Value = 4756.696
$ws.Range("I132").Value = 3082.75
$ws.Range("J132").Value = 5649.467
$ws.Range("K132").Value = 9248.25
$ws.Range("L132").Value = 16948.401
$ws.Range("M132").Value = -6718.25
$ws.Range("N132").Value = -22008.401

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49810
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 49810
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 49810
$ws.Range("N123").Value = -59610

$ws.Range("H126").Value = 35496.83
$ws.Range("I126").Value = 53416.316
$ws.Range("J126").Value = 1449.8
$ws.Range("K126").Value = 160248.948
$ws.Range("L126").Value = 4349.4
$ws.Range("M126").Value = -157778.948
$ws.Range("N126").Value = -9289.4

$ws.Range("H132").Value = 26319438
$ws.Range("I132").Value = 35717464
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 107152392
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -107149862
$ws.Range("N132").Value = -19964
